# Update the Lpl -> Lrp1 sheet with the newly recomputed TPM-based
# expression values (ligand/receptor averages & totals, their derived
# specificities, and the resulting edge weights/specificities for every
# Sending-cluster x Target-cluster combination in rows 2-10).
#
# Columns:
#   G/H  = Ligand (Lpl) average / total expression value for the sending cluster
#   I/J  = Ligand derived specificity (avg / total) for the sending cluster
#   M/N  = Receptor (Lrp1) average / total expression value for the target cluster
#   O/P  = Receptor derived specificity (avg / total) for the target cluster
#   Q/R  = Edge average / total expression weight (ligand * receptor)
#   S/T  = Edge average / total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ G=99.14059966666667;  H=297.421799;         I=0.3911422343348016; J=0.3911422343348016;
             M=3.456265333333333; N=10.368796;           O=0.009841535807677501; P=0.0098415358076775;
             Q=342.6562177537782; R=3083.905959784004;   S=0.003849440305100934; T=0.003849440305100934 }
    3  = @{ G=99.14059966666667;  H=297.421799;         I=0.3911422343348016; J=0.3911422343348016;
             O=0.8587907398420774; P=0.8587907398420773;
             Q=29900.81959836915; R=269107.3763853224;  S=0.3359093288078674; T=0.3359093288078674 }
    4  = @{ G=99.14059966666667;  H=297.421799;         I=0.3911422343348016; J=0.3911422343348016;
             O=0.1313677243502452; P=0.1313677243502452;
             Q=4573.876317724717; R=41164.88685952246;  S=0.05138346522183322; T=0.05138346522183321 }
    5  = @{ I=0.4928190063160421; J=0.4928190063160421;
             M=3.456265333333333; N=10.368796;           O=0.009841535807677501; P=0.0098415358076775;
             Q=431.7291305261777; R=3885.562174735599;   S=0.004850095897363373; T=0.004850095897363372 }
    6  = @{ I=0.4928190063160421; J=0.4928190063160421;
             O=0.8587907398420774; P=0.8587907398420773;
             S=0.4232283990423912; T=0.4232283990423912 }
    7  = @{ I=0.4928190063160421; J=0.4928190063160421;
             O=0.1313677243502452; P=0.1313677243502452;
             S=0.06474051137628757; T=0.06474051137628756 }
    8  = @{ G=29.411685; H=88.235055;                    I=0.1160387593491562; J=0.1160387593491562;
             M=3.456265333333333; N=10.368796;           O=0.009841535807677501; P=0.0098415358076775;
             Q=101.65458726042;   R=914.8912853437801;   S=0.001141999605213193; T=0.001141999605213193 }
    9  = @{ G=29.411685; H=88.235055;                    I=0.1160387593491562; J=0.1160387593491562;
             O=0.8587907398420774; P=0.8587907398420773;
             Q=8870.56856853717;  R=79835.11711683453;   S=0.09965301199181863; T=0.09965301199181861 }
    10 = @{ G=29.411685; H=88.235055;                    I=0.1160387593491562; J=0.1160387593491562;
             O=0.1313677243502452; P=0.1313677243502452;
             Q=1356.915430592355;
             S=0.01524374775212439; T=0.01524374775212438 }
}

foreach ($row in $updates.Keys) {
    $cells = $updates[$row]
    foreach ($col in $cells.Keys) {
        $ws.Range("$col$row").Value2 = $cells[$col]
    }
}
